# 10th - MB for single stock and added new group
#
# The watch sheet keeps one column per refresh date (most recent first in
# column B). This refresh:
#   - drops the three oldest date columns (Jun_24, Jun_22, Jun_19)
#   - inserts two new date columns at the front (Jun_27, Jun_26) with the
#     default "UN" rating for every existing firm row
#   - appends a new group of two firms (Benchmark, Evercore ISI) at the
#     bottom with only the "UN" default rating filled in for the two
#     newest date columns

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember how many data rows currently exist (firms in column A, rows 2..N)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# 1) Drop the three oldest weekly columns (B:D = Jun_24, Jun_22, Jun_19).
#    Everything to the right (Jun_17, Jun_15, Jun_13, Jun_10 ...) shifts left.
$ws.Range("B1:D1").EntireColumn.Delete()

# 2) Insert two fresh blank columns back at the front for the two new
#    refresh dates (Jun_27, Jun_26).
$ws.Range("B1:C1").EntireColumn.Insert()

# 3) Label the two new columns.
$ws.Range("B1").Value = "Jun_27"
$ws.Range("C1").Value = "Jun_26"

# 4) Every existing firm row gets the default "UN" rating for both new
#    columns (no rating data yet for these dates).
for ($row = 2; $row -le $lastRow; $row++) {
    $ws.Cells.Item($row, 2).Value = "UN"
    $ws.Cells.Item($row, 3).Value = "UN"
}

# 5) Add the new group of firms being tracked, each starting at "UN" for
#    the two newest columns only.
$newRow = $lastRow + 1
$ws.Cells.Item($newRow, 1).Value = "Benchmark"
$ws.Cells.Item($newRow, 2).Value = "UN"
$ws.Cells.Item($newRow, 3).Value = "UN"

$newRow = $newRow + 1
$ws.Cells.Item($newRow, 1).Value = "Evercore ISI"
$ws.Cells.Item($newRow, 2).Value = "UN"
$ws.Cells.Item($newRow, 3).Value = "UN"
